# small adjustments to label spacing
#
# Nudge three shapes on slide 1 (PowerPoint's COM object model works in
# points; 1 pt = 12700 EMU). The point values below are chosen so that,
# after the host's internal point->EMU conversion, they land exactly on
# the target EMU (a plain "emu / 12700" can truncate one EMU short because
# of float rounding, so values are nudged by a hair to compensate):
#   - "TextBox 52"             ("Oklab bins" label)  -> reposition (left/top)
#   - "Straight Connector 72"  (divider line)         -> reposition + resize
#   - "TextBox 92"             ("L=0" label)          -> reposition (top only)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Oklab bins" label: off x=43937,y=405906 -> x=79259,y=405672 (ext unchanged)
$oklabLabel = $s.Shapes.Item("TextBox 52")
$oklabLabel.Left = 6.240866141732283   # 79259 EMU
$oklabLabel.Top = 31.94267756535433    # 405672 EMU

# Divider line under the bins row:
#   off x=79259,y=2713521 -> x=150920,y=2713521 (y unchanged)
#   ext cx=11993679,cy=39211 -> cx=11922018,cy=38977
$divider = $s.Shapes.Item("Straight Connector 72")
$divider.Left = 11.883464566929133     # 150920 EMU
$divider.Width = 938.7415748031497     # 11922018 EMU
$divider.Height = 3.069055218110236    # 38977 EMU

# "L=0" label: off x=1655407,y=4357288 -> x=1655407,y=4419433 (x, ext unchanged)
$lLabel = $s.Shapes.Item("TextBox 92")
$lLabel.Top = 347.98686219370074       # 4419433 EMU
